$d = $word.ActiveDocument

function Set-ShapeFillByName($doc, $shapeName, $newRgb) {
    for ($i = 1; $i -le $doc.Shapes.Count; $i++) {
        $s = $doc.Shapes.Item($i)
        if ($s.Name -eq $shapeName) {
            $s.Fill.ForeColor.RGB = $newRgb
        }
    }
}

# Rectangle 264: 02fa47 -> b206c6
Set-ShapeFillByName $d "Rectangle 264" 12977842

# Rectangle 258: c4563d -> f02734
Set-ShapeFillByName $d "Rectangle 258" 3418096

# Rectangle 259: 19cc27 -> 00d2b6
Set-ShapeFillByName $d "Rectangle 259" 11981312

# Rectangle 262: 30bbc1 -> 392ab3
Set-ShapeFillByName $d "Rectangle 262" 11741753
